$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 15
$ws1.Range("F4").Value = 51
$ws1.Range("F5").Value = 26
$ws1.Range("F7").Value = 3471
$ws1.Range("F9").Value = 4131
$ws1.Range("F11").Value = 1035
$ws1.Range("F12").Value = 47

# Sheet "全部类型" (worksheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 15
$ws4.Range("F4").Value = 51
$ws4.Range("F5").Value = 26
$ws4.Range("F8").Value = 3471
$ws4.Range("F10").Value = 4131
$ws4.Range("F12").Value = 1035
$ws4.Range("F13").Value = 47
